# TC09_CDS_Filter_Study-Discovery_CCS.xlsx - "CDS Study filter fixes"
#
# The ParticipantsTab Cypher query (row 2, column B on the "startup" sheet)
# is replaced with a corrected/expanded version of the query (adds
# diagnosis/file/genomic_info traversal, sorts the collected sample ids,
# and reformats the RETURN/ORDER BY/LIMIT clauses). The row grows taller to
# fit the extra lines, and the active selection moves from A2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Discovery of Colorectal Cancer Susceptibility Genes in High-Risk Families"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

# Row 2 = ParticipantsTab: A2 "ParticipantsTab" label stays, B2 gets the
# new query text (replacing the old, shorter one).
$ws.Range("B2").Value = $newParticipantsQuery

# The new query text is taller (more lines), so the row grows to fit it.
$ws.Rows.Item(2).RowHeight = 299.25

# Selection moves onto the (updated) query cell.
$ws.Range("B2").Select()
